$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text in the Price column stays text (avoid float drift / reformatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.639.63"
$ws.Range("E2").Value = "  +6.30%  "
$ws.Range("D3").Value = "1.939.57"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "251.65"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "0.692"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "48.44"
$ws.Range("E8").Value = "  +12.27%  "
$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  +6.89%  "
$ws.Range("D10").Value = "58.77"
$ws.Range("E10").Value = "  +7.68%  "
$ws.Range("D11").Value = "0.0770"
$ws.Range("E11").Value = "  +3.35%  "
$ws.Range("D12").Value = "0.101"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "15.65"
$ws.Range("E13").Value = "  +12.68%  "
$ws.Range("D14").Value = "0.838"
$ws.Range("E14").Value = "  +8.55%  "
$ws.Range("D15").Value = "2.213.99"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D16").Value = "5.18"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").Value = "1.935.15"
$ws.Range("E17").Value = "  +1.98%  "
$ws.Range("D18").Value = "37.519.99"
$ws.Range("E18").Value = "  +5.88%  "
$ws.Range("D19").Value = "75.45"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").Value = "0.0₃0864"
$ws.Range("E20").Value = "  +4.68%  "
$ws.Range("D21").Value = "13.83"
$ws.Range("E21").Value = "  +7.78%  "
$ws.Range("D22").Value = "253.77"
$ws.Range("E22").Value = "  +3.53%  "
$ws.Range("D23").Value = "5.23"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "2.51"
$ws.Range("E25").Value = "  -5.53%  "
$ws.Range("D26").Value = "169.12"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("D27").Value = "2.14"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "8.95"
$ws.Range("E28").Value = "  +4.18%  "
$ws.Range("D29").Value = "18.90"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").Value = "0.129"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").Value = "4.61"
$ws.Range("E31").Value = "  +7.20%  "
$ws.Range("D32").Value = "0.0618"
$ws.Range("E32").Value = "  +3.69%  "
$ws.Range("D33").Value = "0.0920"
$ws.Range("E33").Value = "  +27.14%  "
$ws.Range("D34").Value = "4.36"
$ws.Range("E34").Value = "  +3.49%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "1.88"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "18.97"
$ws.Range("E37").Value = "  +37.48%  "
$ws.Range("D38").Value = "0.906"
$ws.Range("E38").Value = "  +5.65%  "
$ws.Range("D39").Value = "1.45"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").Value = "1.98"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").Value = "106.25"
$ws.Range("E41").Value = "  +7.98%  "
$ws.Range("D42").Value = "0.0228"
$ws.Range("E42").Value = "  +2.89%  "
$ws.Range("D43").Value = "17.64"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").Value = "2.91"
$ws.Range("E44").Value = "  +21.56%  "
$ws.Range("D45").Value = "1.12"
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("D46").Value = "1.353.19"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("D47").Value = "2.43"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").Value = "0.0841"
$ws.Range("E48").Value = "  +4.07%  "
$ws.Range("D49").Value = "2.82"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").Value = "6.46"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").Value = "3.83"
$ws.Range("E51").Value = "  +14.97%  "
